$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 2700
$ws.Range("I20").Value = 2700
$ws.Range("K20").Value = 2700
$ws.Range("M20").Value = -2470
$ws.Range("H35").Value = 2700
$ws.Range("I35").Value = 2700
$ws.Range("K35").Value = 2700
$ws.Range("M35").Value = -2321
$ws.Range("H74").Value = 3286.625
$ws.Range("I74").Value = 2258.8
$ws.Range("K74").Value = 2258.8
$ws.Range("M74").Value = -1322.8
$ws.Range("H77").Value = 3286.625
$ws.Range("I77").Value = 2258.8
$ws.Range("K77").Value = 11294
$ws.Range("M77").Value = -6614
$ws.Range("H98").Value = 1422.6562
$ws.Range("I98").Value = 1040.9565
$ws.Range("K98").Value = 1040.9565
$ws.Range("M98").Value = 457.0435
$ws.Range("H100").Value = 3140.1428
$ws.Range("J100").Value = 4993.3335
$ws.Range("L100").Value = 4993.3335
$ws.Range("N100").Value = -6075.3335
$ws.Range("H107").Value = 1112.9445
$ws.Range("I107").Value = 622.8461
$ws.Range("K107").Value = 622.8461
$ws.Range("M107").Value = 1297.1539
$ws.Range("H112").Value = 2108.9
$ws.Range("J112").Value = 2108.9
$ws.Range("L112").Value = 6326.700000000001
$ws.Range("N112").Value = -8542.700000000001
$ws.Range("H122").Value = 1422.6562
$ws.Range("I122").Value = 1040.9565
$ws.Range("K122").Value = 3122.8695
$ws.Range("M122").Value = -672.8694999999998
$ws.Range("H131").Value = 1699.8948
$ws.Range("J131").Value = 2377.4
$ws.Range("L131").Value = 7132.200000000001
$ws.Range("N131").Value = -17212.2
$ws.Range("H132").Value = 939.7406999999999
$ws.Range("I132").Value = 885.26086
$ws.Range("K132").Value = 2655.78258
$ws.Range("M132").Value = -125.7825800000001
$ws.Range("H138").Value = 1800.0656
$ws.Range("I138").Value = 1260.1538
$ws.Range("J138").Value = 2201.1428
$ws.Range("K138").Value = 3780.4614
$ws.Range("L138").Value = 6603.428400000001
$ws.Range("M138").Value = 1359.5386
$ws.Range("N138").Value = -16883.4284
$ws.Range("H139").Value = 69895.664
$ws.Range("J139").Value = 69895.664
$ws.Range("L139").Value = 69895.664
$ws.Range("N139").Value = -80175.664

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5084.645
$ws.Range("I32").Value = 3033.8655
$ws.Range("J32").Value = 15748.7
$ws.Range("K32").Value = 3033.8655
$ws.Range("L32").Value = 15748.7
$ws.Range("M32").Value = -2746.8655
$ws.Range("N32").Value = -16322.7
$ws.Range("H45").Value = 10001599
$ws.Range("I45").Value = 45000550
$ws.Range("J45").Value = 1899
$ws.Range("K45").Value = 45000550
$ws.Range("L45").Value = 1899
$ws.Range("M45").Value = -45000173
$ws.Range("N45").Value = -2653
$ws.Range("H61").Value = 6765.2383
$ws.Range("I61").Value = 7836.4614
$ws.Range("K61").Value = 7836.4614
$ws.Range("M61").Value = -7624.4614
$ws.Range("H110").Value = 236.18182
$ws.Range("I110").Value = 236.18182
$ws.Range("K110").Value = 236.18182
$ws.Range("M110").Value = 1808.81818
$ws.Range("H132").Value = 1647.9615
$ws.Range("I132").Value = 1175.579
$ws.Range("K132").Value = 3526.737
$ws.Range("M132").Value = -996.7370000000001
$ws.Range("H136").Value = 6765.2383
$ws.Range("I136").Value = 7836.4614
$ws.Range("K136").Value = 23509.3842
$ws.Range("M136").Value = -20959.3842

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2310.353
$ws.Range("I105").Value = 1951.7333
$ws.Range("J105").Value = 5000
$ws.Range("K105").Value = 1951.7333
$ws.Range("L105").Value = 5000
$ws.Range("M105").Value = -204.7333000000001
$ws.Range("N105").Value = -8494
$ws.Range("H134").Value = 10756.917
$ws.Range("I134").Value = 12974.889
$ws.Range("J134").Value = 4103
$ws.Range("K134").Value = 38924.667
$ws.Range("L134").Value = 12309
$ws.Range("M134").Value = -36389.667
$ws.Range("N134").Value = -17379

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 799.2
$ws.Range("I16").Value = 799.2
$ws.Range("K16").Value = 799.2
$ws.Range("M16").Value = -512.2
$ws.Range("H31").Value = 2381.7646
$ws.Range("I31").Value = 2397.5
$ws.Range("J31").Value = 2376.923
$ws.Range("K31").Value = 2397.5
$ws.Range("L31").Value = 2376.923
$ws.Range("M31").Value = -2102.5
$ws.Range("N31").Value = -2966.923
$ws.Range("H34").Value = 2381.7646
$ws.Range("I34").Value = 2397.5
$ws.Range("J34").Value = 2376.923
$ws.Range("K34").Value = 2397.5
$ws.Range("L34").Value = 2376.923
$ws.Range("M34").Value = -2195.5
$ws.Range("N34").Value = -2780.923
$ws.Range("H99").Value = 2034.3636
$ws.Range("I99").Value = 1986.4445
$ws.Range("K99").Value = 1986.4445
$ws.Range("M99").Value = -488.4445000000001
$ws.Range("H113").Value = 799.2
$ws.Range("I113").Value = 799.2
$ws.Range("K113").Value = 799.2
$ws.Range("M113").Value = 1370.8
$ws.Range("H126").Value = 2034.3636
$ws.Range("I126").Value = 1986.4445
$ws.Range("K126").Value = 5959.333500000001
$ws.Range("M126").Value = -3489.333500000001
$ws.Range("H132").Value = 3757.4443
$ws.Range("J132").Value = 4834.6665
$ws.Range("L132").Value = 14503.9995
$ws.Range("N132").Value = -19563.9995
$ws.Range("H134").Value = 3796.111
$ws.Range("I134").Value = 3389.8333
$ws.Range("K134").Value = 10169.4999
$ws.Range("M134").Value = -7634.499899999999
$ws.Range("H138").Value = 104210
$ws.Range("J138").Value = 104210
$ws.Range("L138").Value = 104210
$ws.Range("N138").Value = -114490

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 4000
$ws.Range("J63").Value = 4000
$ws.Range("L63").Value = 12000
$ws.Range("N63").Value = -13498
$ws.Range("H64").Value = 3177.4
$ws.Range("J64").Value = 3177.4
$ws.Range("L64").Value = 9532.200000000001
$ws.Range("N64").Value = -10072.2
$ws.Range("H66").Value = 4000
$ws.Range("J66").Value = 4000
$ws.Range("L66").Value = 36000
$ws.Range("N66").Value = -43488
$ws.Range("H67").Value = 3177.4
$ws.Range("J67").Value = 3177.4
$ws.Range("L67").Value = 9532.200000000001
$ws.Range("N67").Value = -11404.2
$ws.Range("H75").Value = 18878.715
$ws.Range("I75").Value = 270
$ws.Range("J75").Value = 21980.166
$ws.Range("K75").Value = 810
$ws.Range("L75").Value = 65940.49800000001
$ws.Range("M75").Value = 188
$ws.Range("N75").Value = -67936.49800000001
$ws.Range("H78").Value = 18878.715
$ws.Range("I78").Value = 270
$ws.Range("J78").Value = 21980.166
$ws.Range("K78").Value = 2430
$ws.Range("L78").Value = 197821.494
$ws.Range("M78").Value = 2562
$ws.Range("N78").Value = -207805.494
$ws.Range("H131").Value = 14449.98
$ws.Range("I131").Value = 670
$ws.Range("J131").Value = 15622.745
$ws.Range("K131").Value = 2010
$ws.Range("L131").Value = 46868.235
$ws.Range("M131").Value = 3030
$ws.Range("N131").Value = -56948.235
$ws.Range("H132").Value = 903.6667
$ws.Range("I132").Value = 866
$ws.Range("J132").Value = 1205
$ws.Range("K132").Value = 7794
$ws.Range("L132").Value = 10845
$ws.Range("M132").Value = -5264
$ws.Range("N132").Value = -15905

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 95.57143000000001
$ws.Range("I2").Value = 73.8
$ws.Range("J2").Value = 150
$ws.Range("K2").Value = 73.8
$ws.Range("L2").Value = 150
$ws.Range("M2").Value = 39.2
$ws.Range("N2").Value = -376
$ws.Range("H97").Value = 1281
$ws.Range("I97").Value = 1369.5
$ws.Range("K97").Value = 1369.5
$ws.Range("M97").Value = -873.5
$ws.Range("H132").Value = 2567726.2
$ws.Range("J132").Value = 5242.143
$ws.Range("L132").Value = 15726.429
$ws.Range("N132").Value = -20786.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2226.65
$ws.Range("I7").Value = 2112.9443
$ws.Range("K7").Value = 2112.9443
$ws.Range("M7").Value = -2000.9443
$ws.Range("H126").Value = 2226.65
$ws.Range("I126").Value = 2112.9443
$ws.Range("K126").Value = 6338.8329
$ws.Range("M126").Value = -3868.8329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 34910.566
$ws.Range("I122").Value = 44094.5
$ws.Range("J122").Value = 1848.4
$ws.Range("K122").Value = 132283.5
$ws.Range("L122").Value = 5545.200000000001
$ws.Range("M122").Value = -129833.5
$ws.Range("N122").Value = -10445.2
$ws.Range("H132").Value = 1342.3414
$ws.Range("I132").Value = 947.2222
$ws.Range("J132").Value = 4187.2
$ws.Range("K132").Value = 2841.6666
$ws.Range("L132").Value = 12561.6
$ws.Range("M132").Value = -311.6666
$ws.Range("N132").Value = -17621.6
